$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.043333411216736
$ws.Range("B1").Value = 3.603909254074097
$ws.Range("C1").Value = 3.783582210540771
$ws.Range("D1").Value = 3.167818546295166
$ws.Range("E1").Value = 1.276435613632202
